# Fruta / hortaliza, semanal
# Insert a new weekly record (row 134) into the Durazno price sheet, pushing
# the existing rows 134:191 down to 135:192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 134 (shifts 134:191 -> 135:192)
$ws.Rows.Item(134).Insert()

# Populate the new row 134 with the new price record
$ws.Range("A134").Value = 5
$ws.Range("B134").Value = "Macroferia Regional de Talca"
$ws.Range("C134").Value = "Maule"
$ws.Range("D134").Value = 44523
$ws.Range("E134").Value = 7
$ws.Range("F134").Value = "Fruta"
$ws.Range("G134").Value = 100103
$ws.Range("H134").Value = "Frutos de hueso (carozo)"
$ws.Range("I134").Value = 100103004
$ws.Range("J134").Value = "Durazno"
$ws.Range("K134").Value = "Florida King"
$ws.Range("L134").Value = "Especial"
$ws.Range("M134").Value = 15
$ws.Range("N134").Value = 20000
$ws.Range("O134").Value = 20000
$ws.Range("P134").Value = 20000
$ws.Range("Q134").Value = '$/bandeja 15 kilos empedrada'
$ws.Range("R134").Value = "Región de O'Higgins"
$ws.Range("S134").Value = 1333
$ws.Range("T134").Value = 15

# Match the date-number-format style used by the rest of column D
$ws.Range("D134").NumberFormat = "YYYY-MM-DD HH:MM:SS"
